# Applies:
#  1. Update the footer "datetimeFigureOut" date field text, on the slide
#     master and every slide layout, from 6/18/2024 -> 6/24/2024.
#  2. Append a new blank slide (slide id 258) at the end of the deck.

$p = $ppt.ActivePresentation

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shape = $shapes.Item($i)
        if ($shape.HasTextFrame -and $shape.PlaceholderFormat.Type -eq 16) {
            $tr = $shape.TextFrame.TextRange
            if ($tr.Text -eq "6/18/2024") {
                $tr.Text = "6/24/2024"
            }
        }
    }
}

# Slide master footer date placeholder.
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every slide layout's footer date placeholder.
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}

# Add a new blank slide at the end (position 3) -> sldId 258.
$null = $p.Slides.Add($p.Slides.Count + 1, 7)
